$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.863.29'
$ws.Range("E2").Value = '  -0.27%  '
$ws.Range("D3").Value = '1.870.78'
$ws.Range("E3").Value = '  -1.31%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7333'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -5.59%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '242.01'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.09%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9998'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("E8").Value = '  +0.45%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '24.72'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.64%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07093'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.32%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08409'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -8.68%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7518'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.86%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.403'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.24%  '
$ws.Range("D14").Value = '1.854.79'
$ws.Range("E14").Value = '  -3.81%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '92.53'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.26%  '
$ws.Range("D16").Value = '29.865.62'
$ws.Range("E16").Value = '  -0.35%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.043'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.27%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.57'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.97%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '242.83'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.80%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007815'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.87%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9997'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.04%  '
$ws.Range("D22").Value = '2.121.05'
$ws.Range("E22").Value = '  -2.74%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.912'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.31%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1562'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.05%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.313'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.44%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '163.89'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.95%  '
$ws.Range("E28").Value = '  -1.40%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.014'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.81%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.472'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.54%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.644'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.20%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.529'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.50%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.280'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.88%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05333'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.28%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.233'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.43%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7524'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.43%  '
$ws.Range("E37").Value = '  -0.24%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.700'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.34%  '
$ws.Range("E39").Value = '  -0.70%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.752'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.48%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4469'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.00%  '
$ws.Range("D42").Value = '1.104.83'
$ws.Range("E42").Value = '  +1.40%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.080'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.21%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '72.20'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.84%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8624'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.91%  '
$ws.Range("E46").Value = '  +0.07%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '102.84'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.09%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.709'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.32%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.841'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.99%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.062'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.88%  '
$ws.Range("D51").Value = '2.017.59'
$ws.Range("E51").Value = '  -2.62%  '
